$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row 22 (top table) : "RunningLeftKick" ---
$ws.Range("A22").Value = "RunningLeftKick"
$ws.Range("B22").Formula = "=B46"
$ws.Range("C22").Formula = "=F22+B22-E22+5"
$ws.Range("D22").Formula = "=G22+B22-E22"
$ws.Range("E22").Formula = "=E46"
$ws.Range("F22").Formula = "=F46"
$ws.Range("G22").Formula = "=G46"

# --- New row 46 (bottom/source table) : "RunningLeftKick" ---
$ws.Range("A46").Value = "RunningLeftKick"
$ws.Range("B46").Value = 25
$ws.Range("C46").Value = 15
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 65
$ws.Range("F46").Formula = "=C46-B46+E46-5"
$ws.Range("G46").Formula = "=D46-B46+E46"

# --- View state: scroll so row 16 is at top, select E46 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E46").Select()
